$d = $word.ActiveDocument

$replacements = @(
    @{old="40×96=3840"; new="11×71=781"},
    @{old="86×15=1290"; new="43×93=3999"},
    @{old="77×11=847"; new="84×89=7476"},
    @{old="38×25=950"; new="56×21=1176"},
    @{old="85×19=1615"; new="69×26=1794"},
    @{old="44×66=2904"; new="56×71=3976"},
    @{old="55×89=4895"; new="83×93=7719"},
    @{old="99×35=3465"; new="78×57=4446"},
    @{old="35×12=420"; new="64×13=832"},
    @{old="61×34=2074"; new="55×28=1540"},
    @{old="99×73=7227"; new="37×97=3589"},
    @{old="43×51=2193"; new="64×83=5312"},
    @{old="76×29=2204"; new="18×87=1566"},
    @{old="91×39=3549"; new="42×53=2226"},
    @{old="25×76=1900"; new="49×26=1274"},
    @{old="92×43=3956"; new="81×79=6399"},
    @{old="52×18=936"; new="45×37=1665"},
    @{old="70×13=910"; new="45×44=1980"},
    @{old="58×52=3016"; new="45×12=540"},
    @{old="73×93=6789"; new="56×28=1568"},
    @{old="16×35=560"; new="15×79=1185"},
    @{old="73×54=3942"; new="28×62=1736"},
    @{old="17×25=425"; new="24×61=1464"},
    @{old="60×48=2880"; new="78×53=4134"},
    @{old="27×57=1539"; new="54×34=1836"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
